$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 27.5
$ws.Range("E2").Value = 28.06999969482422
$ws.Range("F2").Value = 34.83000183105469
$ws.Range("G2").Value = 26.05999946594238
$ws.Range("H2").Value = 158300823
$ws.Range("I2").Value = "ZS"
$ws.Range("D3").Value = 27.5
$ws.Range("E3").Value = 28.06999969482422
$ws.Range("F3").Value = 34.83000183105469
$ws.Range("G3").Value = 26.05999946594238
$ws.Range("H3").Value = 158300823
$ws.Range("I3").Value = "ZS"
$ws.Range("D4").Value = 27.5
$ws.Range("E4").Value = 28.06999969482422
$ws.Range("F4").Value = 34.83000183105469
$ws.Range("G4").Value = 26.05999946594238
$ws.Range("H4").Value = 158300823
$ws.Range("I4").Value = "ZS"
$ws.Range("D5").Value = 27.5
$ws.Range("E5").Value = 28.06999969482422
$ws.Range("F5").Value = 34.83000183105469
$ws.Range("G5").Value = 26.05999946594238
$ws.Range("H5").Value = 158300823
$ws.Range("I5").Value = "ZS"
$ws.Range("D6").Value = 27.5
$ws.Range("E6").Value = 28.06999969482422
$ws.Range("F6").Value = 34.83000183105469
$ws.Range("G6").Value = 26.05999946594238
$ws.Range("H6").Value = 158300823
$ws.Range("I6").Value = "ZS"
$ws.Range("D8").Value = 35.59999847412109
$ws.Range("E8").Value = 42.79000091552734
$ws.Range("F8").Value = 48.2400016784668
$ws.Range("G8").Value = 35.19100189208984
$ws.Range("H8").Value = 158300823
$ws.Range("I8").Value = "ZS"
$ws.Range("D9").Value = 36.29999923706055
$ws.Range("E9").Value = 39.2599983215332
$ws.Range("F9").Value = 42.47999954223633
$ws.Range("G9").Value = 30.71999931335449
$ws.Range("H9").Value = 158300823
$ws.Range("I9").Value = "ZS"
$ws.Range("D10").Value = 48.04000091552734
$ws.Range("E10").Value = 49.68000030517578
$ws.Range("F10").Value = 51.36000061035156
$ws.Range("G10").Value = 46.5099983215332
$ws.Range("H10").Value = 158300823
$ws.Range("I10").Value = "ZS"
$ws.Range("D11").Value = 68.9800033569336
$ws.Range("E11").Value = 68.62999725341797
$ws.Range("F11").Value = 80.37999725341797
$ws.Range("G11").Value = 63.40499877929688
$ws.Range("H11").Value = 158300823
$ws.Range("I11").Value = "ZS"
$ws.Range("D12").Value = 84.83000183105469
$ws.Range("E12").Value = 68.73999786376953
$ws.Range("F12").Value = 86.41000366210938
$ws.Range("G12").Value = 67.77999877929688
$ws.Range("H12").Value = 158300823
$ws.Range("I12").Value = "ZS"
$ws.Range("D13").Value = 44
$ws.Range("E13").Value = 52.13000106811523
$ws.Range("F13").Value = 52.68000030517578
$ws.Range("G13").Value = 42.29000091552734
$ws.Range("H13").Value = 158300823
$ws.Range("I13").Value = "ZS"
$ws.Range("D14").Value = 56.56000137329102
$ws.Range("E14").Value = 51.9900016784668
$ws.Range("F14").Value = 66.5
$ws.Range("G14").Value = 48.04000091552734
$ws.Range("H14").Value = 158300823
$ws.Range("I14").Value = "ZS"
$ws.Range("D15").Value = 66.2300033569336
$ws.Range("E15").Value = 98.08999633789062
$ws.Range("F15").Value = 98.37999725341795
$ws.Range("G15").Value = 65.51000213623047
$ws.Range("H15").Value = 158300823
$ws.Range("I15").Value = "ZS"
$ws.Range("D16").Value = 131.0500030517578
$ws.Range("E16").Value = 143.3399963378906
$ws.Range("F16").Value = 146.9299926757812
$ws.Range("G16").Value = 116.879997253418
$ws.Range("H16").Value = 158300823
$ws.Range("I16").Value = "ZS"
$ws.Range("D17").Value = 136.75
$ws.Range("E17").Value = 155.75
$ws.Range("F17").Value = 155.75
$ws.Range("G17").Value = 120.3399963378906
$ws.Range("H17").Value = 158300823
$ws.Range("I17").Value = "ZS"
$ws.Range("D18").Value = 202.5200042724609
$ws.Range("E18").Value = 205.0299987792969
$ws.Range("F18").Value = 230.8809967041016
$ws.Range("G18").Value = 190.6000061035156
$ws.Range("H18").Value = 158300823
$ws.Range("I18").Value = "ZS"
$ws.Range("D19").Value = 186.6000061035156
$ws.Range("E19").Value = 194.1999969482422
$ws.Range("F19").Value = 199.7700042724609
$ws.Range("G19").Value = 157.0299987792969
$ws.Range("H19").Value = 158300823
$ws.Range("I19").Value = "ZS"
$ws.Range("D20").Value = 237
$ws.Range("E20").Value = 278.3399963378906
$ws.Range("F20").Value = 279.2900085449219
$ws.Range("G20").Value = 229.5
$ws.Range("H20").Value = 158300823
$ws.Range("I20").Value = "ZS"
$ws.Range("D21").Value = 319.4599914550781
$ws.Range("E21").Value = 346.9700012207031
$ws.Range("F21").Value = 376.1099853515625
$ws.Range("G21").Value = 315.1300048828125
$ws.Range("H21").Value = 158300823
$ws.Range("I21").Value = "ZS"
$ws.Range("D22").Value = 260
$ws.Range("E22").Value = 239.1499938964844
$ws.Range("F22").Value = 290.4700012207031
$ws.Range("G22").Value = 206
$ws.Range("H22").Value = 158300823
$ws.Range("I22").Value = "ZS"
$ws.Range("D23").Value = 201.3099975585937
$ws.Range("E23").Value = 153.0899963378906
$ws.Range("F23").Value = 210.6600036621093
$ws.Range("G23").Value = 125.120002746582
$ws.Range("H23").Value = 158300823
$ws.Range("I23").Value = "ZS"
$ws.Range("D24").Value = 152
$ws.Range("E24").Value = 159.2400054931641
$ws.Range("F24").Value = 188.0599975585937
$ws.Range("G24").Value = 149.0800018310547
$ws.Range("H24").Value = 158300823
$ws.Range("I24").Value = "ZS"
$ws.Range("D25").Value = 159
$ws.Range("E25").Value = 133.4499969482422
$ws.Range("F25").Value = 159.5500030517578
$ws.Range("G25").Value = 114.2099990844727
$ws.Range("H25").Value = 158300823
$ws.Range("I25").Value = "ZS"
$ws.Range("D26").Value = 125.0500030517578
$ws.Range("E26").Value = 131.1499938964844
$ws.Range("F26").Value = 143.6300048828125
$ws.Range("G26").Value = 123.9100036621094
$ws.Range("H26").Value = 158300823
$ws.Range("I26").Value = "ZS"
$ws.Range("D27").Value = 89.4800033569336
$ws.Range("E27").Value = 135.4799957275391
$ws.Range("F27").Value = 137.6499938964844
$ws.Range("G27").Value = 84.93000030517578
$ws.Range("H27").Value = 158300823
$ws.Range("I27").Value = "ZS"
$ws.Range("D28").Value = 159
$ws.Range("E28").Value = 156.0500030517578
$ws.Range("F28").Value = 164.2899932861328
$ws.Range("G28").Value = 131.5899963378906
$ws.Range("H28").Value = 158300823
$ws.Range("I28").Value = "ZS"
$ws.Range("D29").Value = 159.6000061035156
$ws.Range("E29").Value = 197.5299987792969
$ws.Range("F29").Value = 207.781005859375
$ws.Range("G29").Value = 152.6000061035156
$ws.Range("H29").Value = 158300823
$ws.Range("I29").Value = "ZS"
$ws.Range("D30").Value = 236.6100006103516
$ws.Range("E30").Value = 241.9700012207031
$ws.Range("F30").Value = 259.6099853515625
$ws.Range("G30").Value = 202.4400024414062
$ws.Range("H30").Value = 158300823
$ws.Range("I30").Value = "ZS"
$ws.Range("D31").Value = 173
$ws.Range("E31").Value = 169.9600067138672
$ws.Range("F31").Value = 183.6999969482422
$ws.Range("G31").Value = 155.25
$ws.Range("H31").Value = 158300823
$ws.Range("I31").Value = "ZS"
$ws.Range("D32").Value = 178.6199951171875
$ws.Range("E32").Value = 199.979995727539
$ws.Range("F32").Value = 202.6399993896484
$ws.Range("G32").Value = 155.9700012207031
$ws.Range("H32").Value = 158300823
$ws.Range("I32").Value = "ZS"
$ws.Range("D33").Value = 181.1499938964844
$ws.Range("E33").Value = 206.5899963378907
$ws.Range("F33").Value = 217.8399963378907
$ws.Range("G33").Value = 179.6300048828125
$ws.Range("H33").Value = 158300823
$ws.Range("I33").Value = "ZS"
$ws.Range("D34").Value = 198.8800048828125
$ws.Range("E34").Value = 196.229995727539
$ws.Range("F34").Value = 217.1000061035156
$ws.Range("G34").Value = 189.0200042724609
$ws.Range("H34").Value = 158300823
$ws.Range("I34").Value = "ZS"
$ws.Range("D35").Value = 228.5
$ws.Range("E35").Value = 275.7000122070312
$ws.Range("F35").Value = 276.489990234375
$ws.Range("G35").Value = 226.5449981689453
$ws.Range("H35").Value = 158300823
$ws.Range("I35").Value = "ZS"
$ws.Range("D36").Value = 281.4500122070312
$ws.Range("E36").Value = 277.0499877929688
$ws.Range("F36").Value = 291.0700073242188
$ws.Range("G36").Value = 265.8599853515625
$ws.Range("H36").Value = 158300823
$ws.Range("I36").Value = "ZS"
